$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$rng = $ws1.Range("A1")
$rng.Font.ThemeColor = 1
Write-Output $rng.Font.ThemeColor
